$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "Thu Nov 07 16:45:57 EST 2024"
$ws.Range("B3").Value = "Thu Nov 07 16:46:13 EST 2024"
$ws.Range("B4").Value = "Thu Nov 07 16:46:28 EST 2024"
